$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(0.9836888475726653, 0.9666629224337555, 0.956697222465607, 0.9527592378627503, 0.952112784532801, 0.9566436144769455, 0.9777171573561247, 1.022905267140629, 1.058449588371786, 1.075126889773486, 1.081514946874563, 1.080135935190299, 1.075650983419933, 1.072913282252955, 1.057369882037477, 1.047964390401802, 1.042602436187934, 1.040795196451256, 1.048960672588436, 1.076966351326575, 1.095693428431559, 1.085659769318198, 1.048510112007705, 1.010268208455443)
$colD = @(0.05153266548520641, 0.05135146586769856, 0.05123492601603807, 0.05118611249832838, 0.0511779273418167, 0.05123427304789807, 0.05147128762207132, 0.05189392209198651, 0.05217844395249926, 0.05230217971328344, 0.05234821142714097, 0.05233833442117941, 0.05230598332124714, 0.05228605982800971, 0.05217024246800861, 0.05209772979390692, 0.05205548672577365, 0.05204109211290842, 0.05210550437408656, 0.05231550803778262, 0.05244795032794869, 0.05237770500065153, 0.05210199121325232, 0.05178413162492745)
$colE = @(0.3485376677749841, 0.3494438080191173, 0.3500366861709487, 0.3502874814278365, 0.3503296813908046, 0.3500400312445628, 0.3488425364638728, 0.3467832930061361, 0.3454457260728079, 0.3448751306452174, 0.3446644935549195, 0.3447096164670223, 0.344857692540125, 0.3449491008940102, 0.3454837769582753, 0.3458214752823228, 0.3460192756111224, 0.3460868600552853, 0.3457851578147171, 0.3448140515880835, 0.3442110520673536, 0.3445299895462757, 0.3458015655633853, 0.3473095171989353)
$colF = @(1.358309705618936, 1.340980092912403, 1.331034923750934, 1.327156793782251, 1.326523372481475, 1.330981915262598, 1.352190003998459, 1.399309704366047, 1.437324835458057, 1.455361994108173, 1.462299508231325, 1.460800616995499, 1.455930596939893, 1.452961541839088, 1.436161056152969, 1.426045217025575, 1.420296856120629, 1.418362574793349, 1.427114818486487, 1.457358127491389, 1.477749039854217, 1.46680873678217, 1.426631041826539, 1.385967930351626)
$colG = @(1.248498679635105, 1.231150063216816, 1.221156276963427, 1.217248931866621, 1.216610088185334, 1.221102912691123, 1.242380027216512, 1.289346779375734, 1.327080225509434, 1.344953833601693, 1.351824437482549, 1.350340176102009, 1.345517030781195, 1.342576043033972, 1.325926431886131, 1.315894161666478, 1.31019052752734, 1.308270812786105, 1.316955211905821, 1.346930927252998, 1.367117993866145, 1.356289098886805, 1.316475311599817, 1.276076579724815)
$colH = @(1.165324874991541, 1.161908656791923, 1.160292207186899, 1.159754429404742, 1.159672434957088, 1.160284464915719, 1.164047076273008, 1.175246603054603, 1.185811575852682, 1.191126965139603, 1.193213101435219, 1.192760552846721, 1.191297123094273, 1.190410279890159, 1.185474457889995, 1.182576994231084, 1.180958385946553, 1.180418583876815, 1.18288047243891, 1.191724977591321, 1.197932735568116, 1.194580410487191, 1.182743122967082, 1.171807014876862)
$colK = @(0.3589642046039501, 0.3128378069403368, 0.2844377724248375, 0.2728449572488785, 0.2709188003141776, 0.2842815067619142, 0.3430761084274252, 0.4577526047998504, 0.5416380266496787, 0.5797226925699306, 0.594133590481448, 0.5910304341271342, 0.5809085057094876, 0.5747071042875973, 0.5391475799065404, 0.5173136074842546, 0.5047482096130125, 0.5004925704924119, 0.5196386041863548, 0.5838818592468442, 0.6258047297558562, 0.6034355818705706, 0.5185875129093631, 0.4267948550564427)
$colL = @(0.1130745574415073, 0.107779728634199, 0.1045943035808605, 0.103312748721109, 0.1031009469105157, 0.1045769531058838, 0.1112352970878447, 0.1248124601969636, 0.1351052854706154, 0.1398569505229688, 0.1416662541148952, 0.1412761460975105, 0.1400056037495432, 0.1392286549432953, 0.1347961471002463, 0.1320947027881374, 0.1305474381372278, 0.1300246841549892, 0.1323816000495697, 0.1403785229958459, 0.145662964106009, 0.1428372619811, 0.1322518756311695, 0.1210837328456051)
$colM = @(0.1976515154339644, 0.1927531806438658, 0.1898510190654612, 0.1886949443965911, 0.1885045864193664, 0.1898353201282319, 0.1959407244030338, 0.2087479269541177, 0.2186646610085603, 0.2232859734486112, 0.2250517370251259, 0.2246707477757113, 0.2234309281407789, 0.2226735549807799, 0.2183648580915687, 0.2157497751285078, 0.2142560193707581, 0.2137520425276591, 0.2160270823101342, 0.2237946656673842, 0.2289631368159846, 0.2261962384459721, 0.2159016815446648, 0.2051940516405395)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 4).Value = $colD[$i]
    $ws.Cells.Item($row, 5).Value = $colE[$i]
    $ws.Cells.Item($row, 6).Value = $colF[$i]
    $ws.Cells.Item($row, 7).Value = $colG[$i]
    $ws.Cells.Item($row, 8).Value = $colH[$i]
    $ws.Cells.Item($row, 11).Value = $colK[$i]
    $ws.Cells.Item($row, 12).Value = $colL[$i]
    $ws.Cells.Item($row, 13).Value = $colM[$i]
}

Write-Output "applied"
